# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right after "总计" by duplicating the layout
#    of the existing "2022-Q3" sheet (so all styles/margins/page setup match), then
#    overwrite its cells with the new quarter's fund-holding data.
# 2) Insert a new row into the "总计" (total) summary sheet for the 2022-Q4 quarter,
#    above the existing 2022-Q3 row, and re-number the leading index column.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet (keeps identical styling/margins) and place
#    the copy right after "总计"; rename it to "2022-Q4".
# ---------------------------------------------------------------------------
$q3Sheet.Copy([System.Reflection.Missing]::Value, $total)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The source sheet only had 2 data rows (rows 2-3); we need 5, so clone row 3's
# formatting down through row 6.
$newSheet.Range("A3:H3").Copy()
$newSheet.Range("A4:H6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings detail for 2022-Q4
$rows = @(
    @(0, "090001", "大成价值增长混合", "12.94", "66.15", "4.77", "0.6172", 4),
    @(1, "160919", "大成产业升级股票（LOF）", "3.30", "86.93", "4.89", "0.1614", 7),
    @(2, "012051", "申万菱信乐道三年持有期混合", "3.48", "87.51", "3.51", "0.1221", 10),
    @(3, "012210", "申万菱信智能汽车股票A", "2.83", "85.83", "4.18", "0.1183", 10),
    @(4, "012211", "申万菱信智能汽车股票C", "1.03", "85.83", "4.18", "0.0431", 10)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = "'" + $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = "'" + $row[3]
    $newSheet.Range("E$r").Value = "'" + $row[4]
    $newSheet.Range("F$r").Value = "'" + $row[5]
    $newSheet.Range("G$r").Value = "'" + $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q4 summary row into the "总计" sheet, right above 2022-Q3,
#    and re-number the leading index column (0,1,2,...).
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# A3 (the old A2, now shifted down) already carries the correct index-column
# style (s=2); clone it onto the freshly inserted A2 cell.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.06

for ($i = 3; $i -le 9; $i++) {
    $total.Range("A$i").Value = $i - 2
}

Write-Host "2022-Q4 sheet and summary row added"
